$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the AccessKey value in D2 with the new test access key
$ws.Range("D2").Value = "7fe67bf08c80ded756e598d6f8fedaea"

# Move selection/active cell to A2
$ws.Range("A2").Select()
